$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'300.91"
$ws.Range("E2").Value = "'-1.02%"
$ws.Range("D3").Value = "'36.55"
$ws.Range("E3").Value = "'2.34%"
$ws.Range("D4").Value = "'4.992"
$ws.Range("E4").Value = "'-1.84%"
$ws.Range("D5").Value = "'0.07683"
$ws.Range("E5").Value = "'-1.52%"
$ws.Range("D6").Value = "'2.057"
$ws.Range("E6").Value = "'-9.20%"
$ws.Range("D7").Value = "'7.918"
$ws.Range("E7").Value = "'-2.03%"
$ws.Range("D8").Value = "'4.026"
$ws.Range("E8").Value = "'-0.52%"
$ws.Range("D9").Value = "'0.9184"
$ws.Range("E9").Value = "'-1.13%"
$ws.Range("D10").Value = "'0.09662"
$ws.Range("E10").Value = "'7.52%"
$ws.Range("D11").Value = "'0.1861"
$ws.Range("E11").Value = "'1.40%"
$ws.Range("D12").Value = "'0.08494"
$ws.Range("E12").Value = "'-0.95%"
$ws.Range("D13").Value = "'0.03524"
$ws.Range("E13").Value = "'-6.61%"
$ws.Range("D14").Value = "'0.09954"
$ws.Range("E14").Value = "'0.19%"
$ws.Range("D15").Value = "'0.001480"
$ws.Range("E15").Value = "'0.33%"
$ws.Range("D16").Value = "'0.005645"
$ws.Range("E16").Value = "'-0.85%"
$ws.Range("D18").Value = "'2.420"
$ws.Range("E18").Value = "'10.80%"
$ws.Range("D19").Value = "'0.3383"
$ws.Range("E19").Value = "'-2.31%"
$ws.Range("D20").Value = "'0.1327"
$ws.Range("E20").Value = "'0.36%"
$ws.Range("D21").Value = "'4.759"
$ws.Range("E21").Value = "'3.76%"
$ws.Range("E22").Value = "'-1.61%"
$ws.Range("D23").Value = "'0.04602"
$ws.Range("E23").Value = "'-1.60%"
$ws.Range("D24").Value = "'0.005086"
$ws.Range("E24").Value = "'12.28%"
$ws.Range("E25").Value = "'-0.08%"
$ws.Range("D26").Value = "'0.0001401"
$ws.Range("E26").Value = "'7.51%"
$ws.Range("D39").Value = "'0.01755"
$ws.Range("E39").Value = "'-0.69%"
$ws.Range("D40").Value = "'0.04604"
$ws.Range("E40").Value = "'-2.91%"
$ws.Range("D41").Value = "'0.007497"
$ws.Range("E41").Value = "'-5.47%"
$ws.Range("D42").Value = "'0.1388"
$ws.Range("E42").Value = "'-1.81%"
$ws.Range("D43").Value = "'0.007725"
$ws.Range("E43").Value = "'-3.37%"
$ws.Range("E44").Value = "'-2.77%"
$ws.Range("D45").Value = "'0.01034"
$ws.Range("E45").Value = "'7.37%"
$ws.Range("D46").Value = "'0.00006277"
$ws.Range("E46").Value = "'0.87%"
$ws.Range("E47").Value = "'-0.18%"
$ws.Range("D48").Value = "'0.0005803"
$ws.Range("E48").Value = "'0.05%"
$ws.Range("D49").Value = "'35.25"
$ws.Range("E49").Value = "'509.16%"
$ws.Range("D50").Value = "'0.002001"
$ws.Range("E50").Value = "'-25.79%"
$ws.Range("D51").Value = "'0.00002102"
$ws.Range("E51").Value = "'-0.18%"
Write-Output "Updated symbol list"
